# Apply updates to the "Data" worksheet, replacing old sequential
# transfer/upload id numbers with the new ones.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rows that have a value in column K (60000347 -> 60000373)
$kRows = @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29)
foreach ($r in $kRows) {
    $ws.Range("K$r").Value = 60000373
}

# Rows that have a value in column L (257 -> 271)
$lRows = @(7, 8, 10, 12, 13, 15, 17, 18, 20, 22, 23, 25, 27, 28)
foreach ($r in $lRows) {
    $ws.Range("L$r").Value = 271
}

# Rows that have a value in column N (60000348 -> 60000374)
$nRows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)
foreach ($r in $nRows) {
    $ws.Range("N$r").Value = 60000374
}

# Rows that have a value in column O (258 -> 272)
$oRows = @(8, 13, 18, 23, 28)
foreach ($r in $oRows) {
    $ws.Range("O$r").Value = 272
}
